$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44253
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 25000
$ws.Range("L2").Value = 26000
$ws.Range("M2").Value = 25500
$ws.Range("P2").Value = 1020
$ws.Range("D3").Value = 44230
$ws.Range("H3").Value = 'Magnum'
$ws.Range("K3").Value = 22000
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = 23000
$ws.Range("P3").Value = 920
$ws.Range("D4").Value = 44321
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24500
$ws.Range("N4").Value = '$/saco 25 kilos'
$ws.Range("O4").Value = 'Región del Maule'
$ws.Range("P4").Value = 980
$ws.Range("D5").Value = 44265
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 21000
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 840
$ws.Range("D6").Value = 44203
$ws.Range("N6").Value = '$/saco 25 kilos'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("D7").Value = 44475
$ws.Range("K7").Value = 44000
$ws.Range("L7").Value = 45000
$ws.Range("M7").Value = 44500
$ws.Range("O7").Value = 'Perú'
$ws.Range("P7").Value = 1780
$ws.Range("D8").Value = 44363
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 26000
$ws.Range("M8").Value = 25500
$ws.Range("O8").Value = 'Perú'
$ws.Range("P8").Value = 1020
$ws.Range("D9").Value = 44441
$ws.Range("H9").Value = 'Magnum'
$ws.Range("K9").Value = 28000
$ws.Range("L9").Value = 29000
$ws.Range("M9").Value = 28500
$ws.Range("N9").Value = '$/malla 25 kilos'
$ws.Range("O9").Value = 'Perú'
$ws.Range("P9").Value = 1140
$ws.Range("D10").Value = 44237
$ws.Range("H10").Value = 'Sin especificar'
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 21000
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 840
$ws.Range("D11").Value = 44489
$ws.Range("K11").Value = 40000
$ws.Range("L11").Value = 42000
$ws.Range("M11").Value = 41000
$ws.Range("P11").Value = 1640
$ws.Range("D12").Value = 44461
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("K12").Value = 33000
$ws.Range("L12").Value = 34000
$ws.Range("M12").Value = 33500
$ws.Range("N12").Value = '$/malla 25 kilos'
$ws.Range("O12").Value = 'Región de Arica y Parinacota'
$ws.Range("P12").Value = 1340
$ws.Range("D13").Value = 44167
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("K13").Value = 18000
$ws.Range("L13").Value = 19000
$ws.Range("M13").Value = 18500
$ws.Range("P13").Value = 740
$ws.Range("D14").Value = 44447
$ws.Range("K14").Value = 37000
$ws.Range("L14").Value = 38000
$ws.Range("M14").Value = 37500
$ws.Range("N14").Value = '$/malla 25 kilos'
$ws.Range("O14").Value = 'Perú'
$ws.Range("P14").Value = 1500
$ws.Range("D15").Value = 44435
$ws.Range("K15").Value = 25000
$ws.Range("L15").Value = 26000
$ws.Range("M15").Value = 25500
$ws.Range("N15").Value = '$/malla 25 kilos'
$ws.Range("O15").Value = 'Perú'
$ws.Range("P15").Value = 1020
$ws.Range("D16").Value = 44468
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("K16").Value = 31000
$ws.Range("L16").Value = 32000
$ws.Range("M16").Value = 31500
$ws.Range("O16").Value = 'Región de Arica y Parinacota'
$ws.Range("P16").Value = 1260
$ws.Range("D17").Value = 44335
$ws.Range("K17").Value = 35000
$ws.Range("L17").Value = 36000
$ws.Range("M17").Value = 35500
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 1420
$ws.Range("D18").Value = 44188
$ws.Range("K18").Value = 38000
$ws.Range("L18").Value = 40000
$ws.Range("M18").Value = 39000
$ws.Range("O18").Value = 'Región Metropolitana'
$ws.Range("P18").Value = 1560
$ws.Range("D19").Value = 44160
$ws.Range("K19").Value = 28000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 29000
$ws.Range("O19").Value = 'Región de O''Higgins'
$ws.Range("P19").Value = 1160
$ws.Range("D20").Value = 44272
$ws.Range("K20").Value = 22000
$ws.Range("L20").Value = 24000
$ws.Range("M20").Value = 23000
$ws.Range("P20").Value = 920
$ws.Range("D21").Value = 44294
$ws.Range("K21").Value = 24000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 24500
$ws.Range("N21").Value = '$/saco 25 kilos'
$ws.Range("O21").Value = 'Región del Maule'
$ws.Range("P21").Value = 980
$ws.Range("D22").Value = 44399
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 22000
$ws.Range("M22").Value = 21000
$ws.Range("P22").Value = 840
$ws.Range("D23").Value = 44279
$ws.Range("K23").Value = 28000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 29000
$ws.Range("N23").Value = '$/saco 25 kilos'
$ws.Range("O23").Value = 'Región del Maule'
$ws.Range("P23").Value = 1160
$ws.Range("D24").Value = 44323
$ws.Range("H24").Value = 'Magnum'
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 21000
$ws.Range("O24").Value = 'Perú'
$ws.Range("P24").Value = 840
$ws.Range("D26").Value = 44342
$ws.Range("K26").Value = 28000
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = 29000
$ws.Range("O26").Value = 'Región Metropolitana'
$ws.Range("P26").Value = 1160
$ws.Range("D27").Value = 44244
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17000
$ws.Range("P27").Value = 680
$ws.Range("D28").Value = 44433
$ws.Range("H28").Value = 'Magnum'
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = 25500
$ws.Range("O28").Value = 'Perú'
$ws.Range("P28").Value = 1020
